$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.84"
$ws.Range("E2").Value = "'0.76%"
$ws.Range("D3").Value = "'26.82"
$ws.Range("E3").Value = "'-4.23%"
$ws.Range("D4").Value = "'4.741"
$ws.Range("E4").Value = "'-9.15%"
$ws.Range("D5").Value = "'0.05928"
$ws.Range("E5").Value = "'1.06%"
$ws.Range("D6").Value = "'6.662"
$ws.Range("E6").Value = "'-0.77%"
$ws.Range("E7").Value = "'0.40%"
$ws.Range("D8").Value = "'0.9436"
$ws.Range("E8").Value = "'-3.65%"
$ws.Range("D9").Value = "'0.1402"
$ws.Range("E9").Value = "'-0.58%"
$ws.Range("D10").Value = "'0.03811"
$ws.Range("E10").Value = "'9.57%"
$ws.Range("D11").Value = "'0.07105"
$ws.Range("E11").Value = "'-1.01%"
$ws.Range("E12").Value = "'-0.29%"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("D14").Value = "'0.001555"
$ws.Range("E14").Value = "'0.57%"
$ws.Range("D15").Value = "'0.0006037"
$ws.Range("E15").Value = "'-94.29%"
$ws.Range("D16").Value = "'0.006056"
$ws.Range("E16").Value = "'4.51%"
$ws.Range("D17").Value = "'3.498"
$ws.Range("E17").Value = "'-0.07%"
$ws.Range("D18").Value = "'3.203"
$ws.Range("E18").Value = "'-0.50%"
$ws.Range("E19").Value = "'-0.15%"
$ws.Range("D20").Value = "'0.3127"
$ws.Range("E20").Value = "'-1.71%"
$ws.Range("E21").Value = "'0.38%"
$ws.Range("D22").Value = "'3.799"
$ws.Range("E22").Value = "'6.99%"
$ws.Range("D23").Value = "'0.04220"
$ws.Range("E23").Value = "'1.30%"
$ws.Range("D25").Value = "'0.001220"
$ws.Range("E25").Value = "'-0.21%"
$ws.Range("D26").Value = "'0.004287"
$ws.Range("E26").Value = "'-10.66%"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("E28").Value = "'1.94%"
$ws.Range("D40").Value = "'0.03824"
$ws.Range("E40").Value = "'0.31%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006195"
$ws.Range("E41").Value = "'6.48%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1101"
$ws.Range("E42").Value = "'0.08%"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("E43").Value = "'-6.94%"
$ws.Range("E44").Value = "'20.72%"
$ws.Range("E45").Value = "'5.10%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("D47").Value = "'0.08848"
$ws.Range("E47").Value = "'-11.50%"
$ws.Range("D48").Value = "'0.002434"
$ws.Range("E48").Value = "'14.13%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E50").Value = "'-0.05%"
